# Adds a new "Exceptional items" column to the Quarterly results sheet,
# inserted right before the existing "P/l before tax" column (i.e. at
# column L), shifting every following column one position to the right.
# The new column header is populated on the two header rows; the data
# rows are left blank for this new metric (no exceptional-items figures
# were available for the historical quarters).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new blank column at L; everything from L..T shifts to M..U,
# and the sheet's UsedRange/dimension grows from A1:T47 to A1:U47.
$ws.Columns("L:L").Insert()

# Populate the new column's two header rows.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"
